$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.933.45'
$ws.Range("E2").Value = '  +1.28%  '
$ws.Range("D3").Value = '2.286.63'
$ws.Range("E3").Value = '  +2.51%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("E6").Value = '  +2.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.65'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +6.77%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.646'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.18'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0971'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '59.04'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.41'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.73%  '
$ws.Range("E14").Value = '  +2.05%  '
$ws.Range("D15").Value = '2.630.41'
$ws.Range("E15").Value = '  +2.62%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.26'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.872'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.62%  '
$ws.Range("D18").Value = '2.290.25'
$ws.Range("E18").Value = '  +3.05%  '
$ws.Range("D19").Value = '42.832.69'
$ws.Range("E19").Value = '  +1.41%  '
$ws.Range("E20").Value = '  +3.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.94%  '
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.95'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.70%  '
$ws.Range("E24").Value = '  +7.37%  '
$ws.Range("E25").Value = '  +0.43%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.50'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.96%  '
$ws.Range("E27").Value = '  -0.20%  '
$ws.Range("E28").Value = '  +0.37%  '
$ws.Range("E29").Value = '  -0.98%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.19'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '166.72'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.04'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.45'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.52%  '
$ws.Range("E34").Value = '  +4.59%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0824'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.19'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +12.72%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.125'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.60'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +11.64%  '
$ws.Range("E39").Value = '  +2.00%  '
$ws.Range("E40").Value = '  -2.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.46'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +15.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.34'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.43%  '
$ws.Range("E43").Value = '  +3.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.219'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +10.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.83'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.76%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.08'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.90'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.00%  '
$ws.Range("E48").Value = '  +2.58%  '
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.17'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '99.59'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.16%  '
